# Prep for a run of all Trajectories [19/06]
# Rebuild Sheet1 data grid to match the target layout (new FL290 column, new
# switch_10/20/30/60 + skip helper columns replacing the old helper columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear out the previous grid (A1:X20) - we rewrite every needed cell below,
# this guarantees stale cells that no longer belong (e.g. the removed
# " switches FL mid cruise" column) are gone too.
$ws.Range("A1:X20").Clear()

# Row 1
$ws.Range("A1").Value = "AC Type"
$ws.Range("B1").Value = "FL170"
$ws.Range("C1").Value = "FL180"
$ws.Range("D1").Value = "FL190"
$ws.Range("E1").Value = "FL210"
$ws.Range("F1").Value = "FL220"
$ws.Range("G1").Value = "FL230"
$ws.Range("H1").Value = "FL240"
$ws.Range("I1").Value = "FL250"
$ws.Range("J1").Value = "FL270"
$ws.Range("K1").Value = "FL290"
$ws.Range("L1").Value = "FL300"
$ws.Range("M1").Value = "FL310"
$ws.Range("N1").Value = "FL320"
$ws.Range("O1").Value = "FL330"
$ws.Range("P1").Value = "FL340"
$ws.Range("Q1").Value = "FL350"
$ws.Range("R1").Value = "FL360"
$ws.Range("S1").Value = "FL370"
$ws.Range("T1").Value = "FL380"
$ws.Range("U1").Value = "FL390"
$ws.Range("V1").Value = "FL400"
$ws.Range("W1").Value = "FL410"
$ws.Range("X1").Value = "switch_10"
$ws.Range("Y1").Value = "switch_20"
$ws.Range("Z1").Value = "switch_30"
$ws.Range("AA1").Value = "switch_60"
$ws.Range("AB1").Value = "skip"
# Row 2
$ws.Range("A2").Value = "A319"
$ws.Range("B2").Value = 0.464
$ws.Range("E2").Value = 0.505
$ws.Range("F2").Value = 0.512
$ws.Range("G2").Value = 0.525
$ws.Range("H2").Value = 0.535
$ws.Range("O2").Value = 0.652
$ws.Range("P2").Value = 0.67
$ws.Range("Q2").Value = 0.684
$ws.Range("R2").Value = 0.702
$ws.Range("S2").Value = 0.719
$ws.Range("T2").Value = 0.73
$ws.Range("U2").Value = 0.743
$ws.Range("X2").Value = 10
$ws.Range("Y2").Value = 20
$ws.Range("Z2").Value = 30
$ws.Range("AA2").Value = 60
$ws.Range("AB2").Value = "KLM59Z"
# Row 3
$ws.Range("A3").Value = "A320"
$ws.Range("C3").Value = 0.551
$ws.Range("G3").Value = 0.616
$ws.Range("H3").Value = 0.628
$ws.Range("L3").Value = 0.716
$ws.Range("M3").Value = 0.737
$ws.Range("N3").Value = 0.757
$ws.Range("O3").Value = 0.771
$ws.Range("P3").Value = 0.786
$ws.Range("Q3").Value = 0.802
$ws.Range("R3").Value = 0.819
$ws.Range("S3").Value = 0.821
$ws.Range("T3").Value = 0.821
$ws.Range("U3").Value = 0.888
$ws.Range("X3").Value = "DLH08W"
$ws.Range("Y3").Value = "DLH156"
$ws.Range("Z3").Value = "QTR022"
$ws.Range("AA3").Value = "GWI2807"
$ws.Range("AB3").Value = "MIBID"
# Row 4
$ws.Range("A4").Value = "A321"
$ws.Range("P4").Value = 0.775
$ws.Range("Q4").Value = 0.791
$ws.Range("X4").Value = "SAS618"
$ws.Range("Y4").Value = "EIN111"
$ws.Range("Z4").Value = "RYR5008"
$ws.Range("AA4").Value = "PGT424"
# Row 5
$ws.Range("A5").Value = "DH8D"
$ws.Range("G5").Value = 0.99
$ws.Range("M5").Value = 0.99
$ws.Range("Y5").Value = "SAS1042"
$ws.Range("AB5").Value = "DH8D - BAD AIRCRAFT"
# Row 6
$ws.Range("A6").Value = "RJ1H"
$ws.Range("L6").Value = 0.64
$ws.Range("Y6").Value = "TAY011"
$ws.Range("AB6").Value = "AUA381"
# Row 7
$ws.Range("A7").Value = "CRJ9"
$ws.Range("O7").Value = 0.759
$ws.Range("P7").Value = 0.78
$ws.Range("R7").Value = 0.813
$ws.Range("S7").Value = 0.84
$ws.Range("Y7").Value = "TFL219"
$ws.Range("AB7").Value = "AUA522D"
# Row 8
$ws.Range("A8").Value = "B734"
$ws.Range("M8").Value = 0.798
$ws.Range("Q8").Value = 0.821
# Row 9
$ws.Range("A9").Value = "B737"
$ws.Range("L9").Value = 0.724
$ws.Range("N9").Value = 0.761
$ws.Range("O9").Value = 0.777
$ws.Range("P9").Value = 0.796
$ws.Range("R9").Value = 0.819
$ws.Range("S9").Value = 0.821
$ws.Range("T9").Value = 0.821
$ws.Range("U9").Value = 0.821
$ws.Range("V9").Value = 0.821
# Row 10
$ws.Range("A10").Value = "B738"
$ws.Range("C10").Value = 0.563
$ws.Range("D10").Value = 0.566
$ws.Range("M10").Value = 0.737
$ws.Range("Q10").Value = 0.801
$ws.Range("R10").Value = 0.819
$ws.Range("S10").Value = 0.821
$ws.Range("T10").Value = 0.821
$ws.Range("U10").Value = 0.821
$ws.Range("V10").Value = 0.865
$ws.Range("W10").Value = 0.866
# Row 11
$ws.Range("A11").Value = "B736"
$ws.Range("T11").Value = 0.821
# Row 12
$ws.Range("A12").Value = "B77W"
$ws.Range("P12").Value = 0.78
# Row 13
$ws.Range("A13").Value = "B752"
$ws.Range("P13").Value = 0.746
$ws.Range("Q13").Value = 0.758
$ws.Range("R13").Value = 0.782
# Row 14
$ws.Range("A14").Value = "B733"
$ws.Range("C14").Value = 0.547
$ws.Range("M14").Value = 0.712
$ws.Range("N14").Value = 0.734
# Row 15
$ws.Range("A15").Value = "A332"
$ws.Range("T15").Value = 0.859
$ws.Range("U15").Value = 0.859
$ws.Range("V15").Value = 0.859
$ws.Range("W15").Value = 0.89
# Row 16
$ws.Range("A16").Value = "F2TH"
$ws.Range("D16").Value = 0.409
# Row 17
$ws.Range("A17").Value = "A310"
$ws.Range("J17").Value = 0.689
# Row 18
$ws.Range("A18").Value = "E145"
$ws.Range("H18").Value = 0.571
$ws.Range("I18").Value = 0.582
$ws.Range("K18").Value = 0.641
# Row 19
$ws.Range("A19").Value = "B77L"
$ws.Range("Q19").Value = 0.801
$ws.Range("R19").Value = 0.813
$ws.Range("S19").Value = 0.84
# Row 20
$ws.Range("A20").Value = "B763"
$ws.Range("T20").Value = 0.861

# H2 keeps its centered alignment style after the clear/rewrite
$ws.Range("H2").HorizontalAlignment = -4108
$ws.Range("H2").VerticalAlignment = -4108

# Column widths for the new helper columns (auto-fit to content)
$ws.Columns("W:Z").ColumnWidth = 9.7109375
$ws.Columns("AA:AA").ColumnWidth = 20.42578125

# Page setup used for the latest print preview
$ws.PageSetup.Orientation = 1

# Restore the active selection/view state
$ws.Range("K19").Select()

Write-Host "Sheet1 rebuilt"
